# add user_id / lessor_id columns to all tables
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "equip" (sheet2): insert "lessor_id" before the old column C (address)
# and append "user_id" after the last column (username).
# Old:  A=equip_id B=lessor_name C=address D=end_time E=contact F=status G=username
# New:  A=equip_id B=lessor_name C=lessor_id D=address E=end_time F=contact G=status H=username I=user_id
# ---------------------------------------------------------------------------
$wsEquip = $wb.Worksheets.Item("equip")
$wsEquip.Range("C1").EntireColumn.Insert()
$wsEquip.Range("C1").Value = "lessor_id"
$wsEquip.Range("I1").Value = "user_id"

# ---------------------------------------------------------------------------
# Sheet "sale_req" (sheet4): append "lessor_id" after the last used column.
# Old: A=student_id B=equip_id C=end_time D=lessor_name
# New: A=student_id B=equip_id C=end_time D=lessor_name E=lessor_id
# ---------------------------------------------------------------------------
$wsSale = $wb.Worksheets.Item("sale_req")
$wsSale.Range("E1").Value = "lessor_id"

# ---------------------------------------------------------------------------
# Sheet "rent_req" (sheet5): insert "user_id" before old column B (detail) and
# insert "lessor_id" right after lessor_name (which, once the first column is
# inserted, has shifted from E to F).
# Old:  A=username B=detail C=contact D=return_time E=lessor_name F=equip_id G=equip_name H=status
# New:  A=username B=user_id C=detail D=contact E=return_time F=lessor_name G=lessor_id H=equip_id I=equip_name J=status
# ---------------------------------------------------------------------------
$wsRentReq = $wb.Worksheets.Item("rent_req")
$wsRentReq.Range("B1").EntireColumn.Insert()
$wsRentReq.Range("B1").Value = "user_id"
$wsRentReq.Range("G1").EntireColumn.Insert()
$wsRentReq.Range("G1").Value = "lessor_id"

# ---------------------------------------------------------------------------
# Sheet "rent_info" (sheet6): insert "lessor_id" right after lessor_name
# (before the old column D, username) and "user_id" right after username
# (before the old column E, rent_time - which has shifted to F by then).
# Old:  A=equip_id B=equip_name C=lessor_name D=username E=rent_time F=status G=return_time H=end_time
# New:  A=equip_id B=equip_name C=lessor_name D=lessor_id E=username F=user_id G=rent_time H=status I=return_time J=end_time
# ---------------------------------------------------------------------------
$wsRentInfo = $wb.Worksheets.Item("rent_info")
$wsRentInfo.Range("D1").EntireColumn.Insert()
$wsRentInfo.Range("D1").Value = "lessor_id"
$wsRentInfo.Range("F1").EntireColumn.Insert()
$wsRentInfo.Range("F1").Value = "user_id"

# ---------------------------------------------------------------------------
# Selections / active-sheet bookkeeping to mirror the author's final view
# state: User sheet shows a block selection, sale_req/rent_req show a single
# cell, and rent_info ends up as the active (tabSelected) sheet.
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("A2:F16").Select() | Out-Null

$wsEquip.Range("F2").Select() | Out-Null

$wsAuth = $wb.Worksheets.Item("auth_req")
$wsAuth.Range("E4").Select() | Out-Null

$wsSale.Range("E1").Select() | Out-Null

$wsRentReq.Range("B1").Select() | Out-Null

$wsRentInfo.Select() | Out-Null
$wsRentInfo.Range("F1").Select() | Out-Null
